$d = $word.ActiveDocument

# Locate the paragraph ending the bibliography entry "Rio de Janeiro: Elsevier
# Editora, 2007." — the three paragraphs immediately following it (a blank
# paragraph, the "Ver no Jupiter..." line, and the "(c) 2020 ..." footer
# line) are removed by this edit, leaving the blank paragraph and the
# page-break paragraph that originally came after them intact.

$anchorText = "Rio de Janeiro: Elsevier Editora, 2007."
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*$anchorText*") {
        $anchorIndex = $i
    }
}

if ($anchorIndex -gt 0) {
    # Delete the three paragraphs right after the anchor paragraph. Removing
    # from the highest index first keeps the lower indices valid.
    $d.Paragraphs.Item($anchorIndex + 3).Range.Delete()
    $d.Paragraphs.Item($anchorIndex + 2).Range.Delete()
    $d.Paragraphs.Item($anchorIndex + 1).Range.Delete()
}
